$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force text format on column D so numeric-looking strings
# (e.g. "551.12") are not auto-converted to numbers when the Value is set.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '63.652.31'
$ws.Range('E2').Value = '  +5.14%  '
$ws.Range('D3').Value = '3.062.43'
$ws.Range('E3').Value = '  +3.72%  '
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').Value = '551.12'
$ws.Range('E5').Value = '  +5.73%  '
$ws.Range('D6').Value = '139.46'
$ws.Range('E6').Value = '  +8.26%  '
$ws.Range('E7').Value = '  -0.19%  '
$ws.Range('D8').Value = '3.050.49'
$ws.Range('E8').Value = '  +3.76%  '
$ws.Range('D9').Value = '0.501'
$ws.Range('E9').Value = '  +4.43%  '
$ws.Range('E10').Value = '  +2.86%  '
$ws.Range('D11').Value = '6.17'
$ws.Range('E11').Value = '  -0.03%  '
$ws.Range('D12').Value = '0.455'
$ws.Range('E12').Value = '  +4.48%  '
$ws.Range('D13').Value = '0.0000226'
$ws.Range('E13').Value = '  +5.18%  '
$ws.Range('D14').Value = '34.91'
$ws.Range('E14').Value = '  +6.69%  '
$ws.Range('D15').Value = '3.572.24'
$ws.Range('E15').Value = '  +3.70%  '
$ws.Range('D16').Value = '63.680.03'
$ws.Range('E16').Value = '  +4.73%  '
$ws.Range('D17').Value = '3.067.16'
$ws.Range('E17').Value = '  +3.36%  '
$ws.Range('E18').Value = '  -0.43%  '
$ws.Range('D19').Value = '6.72'
$ws.Range('E19').Value = '  +5.30%  '
$ws.Range('D20').Value = '482.91'
$ws.Range('E20').Value = '  +6.88%  '
$ws.Range('D21').Value = '13.51'
$ws.Range('E21').Value = '  +4.26%  '
$ws.Range('D22').Value = '0.683'
$ws.Range('E22').Value = '  +2.99%  '
$ws.Range('D23').Value = '7.18'
$ws.Range('E23').Value = '  +6.90%  '
$ws.Range('D24').Value = '81.76'
$ws.Range('E24').Value = '  +5.67%  '
$ws.Range('D25').Value = '12.48'
$ws.Range('E25').Value = '  +8.10%  '
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').Value = '2.75'
$ws.Range('E27').Value = '  +5.97%  '
$ws.Range('E28').Value = '  +5.94%  '
$ws.Range('D29').Value = '1.99'
$ws.Range('E29').Value = '  +10.61%  '
$ws.Range('E30').Value = '  -0.33%  '
$ws.Range('D31').Value = '26.04'
$ws.Range('E31').Value = '  +3.99%  '
$ws.Range('D32').Value = '1.15'
$ws.Range('E32').Value = '  +3.58%  '
$ws.Range('D33').Value = '2.44'
$ws.Range('E33').Value = '  +9.54%  '
$ws.Range('E34').Value = '  +9.78%  '
$ws.Range('D35').Value = '55.42'
$ws.Range('E35').Value = '  +3.06%  '
$ws.Range('E36').Value = '  +5.37%  '
$ws.Range('D37').Value = '467.51'
$ws.Range('E37').Value = '  +4.60%  '
$ws.Range('D38').Value = '3.165.99'
$ws.Range('E38').Value = '  +1.12%  '
$ws.Range('D39').Value = '0.0816'
$ws.Range('E39').Value = '  +6.72%  '
$ws.Range('D40').Value = '0.0396'
$ws.Range('E40').Value = '  +5.21%  '
$ws.Range('D41').Value = '0.119'
$ws.Range('E41').Value = '  +3.29%  '
$ws.Range('D42').Value = '8.23'
$ws.Range('E42').Value = '  +4.74%  '
$ws.Range('D43').Value = '28.40'
$ws.Range('E43').Value = '  +15.20%  '
$ws.Range('E44').Value = '  +8.27%  '
$ws.Range('D45').Value = '0.251'
$ws.Range('E45').Value = '  +4.83%  '
$ws.Range('E46').Value = '  -0.12%  '
$ws.Range('D47').Value = '2.02'
$ws.Range('E47').Value = '  +7.16%  '
$ws.Range('E48').Value = '  +3.35%  '
$ws.Range('D49').Value = '0.0₃0511'
$ws.Range('E49').Value = '  +3.31%  '
$ws.Range('D50').Value = '115.94'
$ws.Range('E50').Value = '  -2.73%  '
$ws.Range('D51').Value = '2.06'
$ws.Range('E51').Value = '  +7.30%  '

# Restore the original (default) cell style for column D now that the
# text values are safely stored, so no stray formatting is introduced.
$ws.Range("D2:D51").Style = "Normal"

